$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.061724
$ws.Range("H2").Value = 0.185172
$ws.Range("I2").Value = 0.09652262708432048
$ws.Range("J2").Value = 0.09652262708432047
$ws.Range("M2").Value = 36.89194233333333
$ws.Range("N2").Value = 110.675827
$ws.Range("O2").Value = 0.3567095043190808
$ws.Range("P2").Value = 0.3567095043190809
$ws.Range("Q2").Value = 2.277118248582667
$ws.Range("R2").Value = 20.494064237244
$ws.Range("S2").Value = 0.03443053846282344
$ws.Range("T2").Value = 0.03443053846282344
$ws.Range("G3").Value = 0.061724
$ws.Range("H3").Value = 0.185172
$ws.Range("I3").Value = 0.09652262708432048
$ws.Range("J3").Value = 0.09652262708432047
$ws.Range("M3").Value = 42.68037399999999
$ws.Range("O3").Value = 0.4126780562577495
$ws.Range("P3").Value = 0.4126780562577496
$ws.Range("Q3").Value = 2.634403404776
$ws.Range("R3").Value = 23.709630642984
$ws.Range("S3").Value = 0.03983277013004898
$ws.Range("T3").Value = 0.03983277013004898
$ws.Range("G4").Value = 0.061724
$ws.Range("H4").Value = 0.185172
$ws.Range("I4").Value = 0.09652262708432048
$ws.Range("J4").Value = 0.09652262708432047
$ws.Range("M4").Value = 23.85061433333334
$ws.Range("N4").Value = 71.55184300000001
$ws.Range("O4").Value = 0.2306124394231696
$ws.Range("P4").Value = 0.2306124394231696
$ws.Range("Q4").Value = 1.472155319110667
$ws.Range("R4").Value = 13.249397871996
$ws.Range("S4").Value = 0.02225931849144805
$ws.Range("T4").Value = 0.02225931849144805
$ws.Range("I5").Value = 0.8735221647273214
$ws.Range("J5").Value = 0.8735221647273215
$ws.Range("M5").Value = 36.89194233333333
$ws.Range("N5").Value = 110.675827
$ws.Range("O5").Value = 0.3567095043190808
$ws.Range("P5").Value = 0.3567095043190809
$ws.Range("Q5").Value = 20.60774060888711
$ws.Range("R5").Value = 185.469665479984
$ws.Range("S5").Value = 0.3115936583916133
$ws.Range("T5").Value = 0.3115936583916133
$ws.Range("I6").Value = 0.8735221647273214
$ws.Range("J6").Value = 0.8735221647273215
$ws.Range("M6").Value = 42.68037399999999
$ws.Range("O6").Value = 0.4126780562577495
$ws.Range("P6").Value = 0.4126780562577496
$ws.Range("S6").Value = 0.3604834290377327
$ws.Range("T6").Value = 0.3604834290377327
$ws.Range("I7").Value = 0.8735221647273214
$ws.Range("J7").Value = 0.8735221647273215
$ws.Range("M7").Value = 23.85061433333334
$ws.Range("N7").Value = 71.55184300000001
$ws.Range("O7").Value = 0.2306124394231696
$ws.Range("P7").Value = 0.2306124394231696
$ws.Range("Q7").Value = 13.32288956496178
$ws.Range("S7").Value = 0.2014450772979754
$ws.Range("T7").Value = 0.2014450772979754
$ws.Range("I8").Value = 0.02995520818835809
$ws.Range("J8").Value = 0.02995520818835809
$ws.Range("M8").Value = 36.89194233333333
$ws.Range("N8").Value = 110.675827
$ws.Range("O8").Value = 0.3567095043190808
$ws.Range("P8").Value = 0.3567095043190809
$ws.Range("Q8").Value = 0.7066897500232221
$ws.Range("R8").Value = 6.360207750209
$ws.Range("S8").Value = 0.01068530746464409
$ws.Range("T8").Value = 0.01068530746464409
$ws.Range("I9").Value = 0.02995520818835809
$ws.Range("J9").Value = 0.02995520818835809
$ws.Range("M9").Value = 42.68037399999999
$ws.Range("O9").Value = 0.4126780562577495
$ws.Range("P9").Value = 0.4126780562577496
$ws.Range("Q9").Value = 0.8175710175526665
$ws.Range("R9").Value = 7.358139157973998
$ws.Range("S9").Value = 0.01236185708996784
$ws.Range("T9").Value = 0.01236185708996784
$ws.Range("I10").Value = 0.02995520818835809
$ws.Range("J10").Value = 0.02995520818835809
$ws.Range("M10").Value = 23.85061433333334
$ws.Range("N10").Value = 71.55184300000001
$ws.Range("O10").Value = 0.2306124394231696
$ws.Range("P10").Value = 0.2306124394231696
$ws.Range("Q10").Value = 0.4568744179645556
$ws.Range("S10").Value = 0.006908043633746163
$ws.Range("T10").Value = 0.006908043633746164
